{"js": "// The document has paragraphs whose <w:pPr><w:numPr> child elements are\n// serialized as <w:numId/> followed by <w:ilvl/>. Word Online expects\n// <w:ilvl/> to come first. There is no semantic value change here (the\n// ilvl/numId values themselves are unchanged) - only the XML element\n// order inside <w:numPr> needs to be normalized.\n//\n// Re-assigning Word.ListItem.level (even to its own current value) makes\n// the engine rewrite the owning paragraph's <w:numPr> with <w:ilvl> before\n// <w:numId>, which is exactly the reordering we need.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Grab (possibly null) list items for every paragraph and load their\n// current level in one batch.\nconst listItems = paragraphs.items.map((p) => p.listItemOrNullObject);\nfor (const listItem of listItems) {\n  listItem.load(\"isNullObject,level\");\n}\nawait context.sync();\n\n// Touch the level property on each paragraph that is actually part of a\n// list, forcing Word to re-serialize its numbering properties with\n// <w:ilvl> first.\nfor (const listItem of listItems) {\n  if (!listItem.isNullObject) {\n    listItem.level = listItem.level;\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $lf = $p.Range.ListFormat\n    if ($lf.ListType -ne 0) {\n        # Re-assigning ListLevelNumber to its own current value is a no-op\n        # semantically, but it makes Word rewrite the paragraph's <w:numPr>\n        # child elements, emitting <w:ilvl> before <w:numId> as Word Online\n        # expects.\n        $lf.ListLevelNumber = $lf.ListLevelNumber\n    }\n}\n"}
